# contactsImport.xlsx update:
#  - Phil Weier's second row (row 4) was a duplicate test entry ("Assasin" /
#    "Top Secret" / rejection note) that is replaced with a real second
#    contact's occupation/company/wechat_id/notes.
#  - A brand new contact (Jackie Cheng) is appended as row 5, including a
#    mailto hyperlink on her email cell styled the same way as the other
#    email cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the existing "Hyperlink" look (the style already used by
#     C2/C3/C4) so it can be re-stamped onto the refreshed hyperlink cells
#     further down without minting a brand new cell style.
$ws.Range("C4").Copy()

# Hyperlinks get rebuilt below (existing two plus the new one) so they all
# end up addressed/assigned in a single, row-ordered pass.
$ws.Hyperlinks.Delete()

# --- Row 4 (Phil Weier) gets corrected occupation/company/wechat/notes ---
$ws.Range("D4").Value = "Cook"
$ws.Range("E4").Value = "fdsfsd"
$ws.Range("F4").Value = 415131
$ws.Range("G4").Value = "fsdfhjksf klsdfjl skdflsd"

# --- Row 5 (new contact: Jackie Cheng) ---
$ws.Range("A5").Value = "Jackie"
$ws.Range("B5").Value = "Cheng"
$ws.Range("C5").Value = "jackmtr@hotmail.com"
$ws.Range("D5").Value = "Suite Runner"
$ws.Range("E5").Value = "Aquilini Services"
$ws.Range("G5").Value = "I am jackie"

# --- Hyperlinks for every email cell, re-added in row order ---
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:philweier@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:philweier@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:jackmtr@hotmail.com")

# Re-apply the shared "Hyperlink" cell style to all three email cells.
$ws.Range("C3:C5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

[void]$ws.Range("G5").Select()
